# Add example about how to use replace command:
# Update the "Howdy." greeting text to use a proper ellipsis character,
# demonstrating a find/replace style edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GREETINGS")

$ws.Range("D5").Value = "Howdy…"

# Reflect the resulting selection left behind on the sheet after the edit.
$ws.Range("D6").Select()
